$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(383).Insert()

$ws.Cells.Item(383, 1).Value = 3
$ws.Cells.Item(383, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(383, 3).Value = "Coquimbo"
$ws.Cells.Item(383, 4).Value = 44516
$ws.Cells.Item(383, 5).Value = 5
$ws.Cells.Item(383, 6).Value = 100112024
$ws.Cells.Item(383, 7).Value = "Choclo"
$ws.Cells.Item(383, 8).Value = "Dulce o Americano"
$ws.Cells.Item(383, 9).Value = "Primera"
$ws.Cells.Item(383, 10).Value = 73
$ws.Cells.Item(383, 11).Value = 41000
$ws.Cells.Item(383, 12).Value = 42000
$ws.Cells.Item(383, 13).Value = 41521
$ws.Cells.Item(383, 14).Value = "`$/malla 70 unidades"
$ws.Cells.Item(383, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(383, 16).Value = 593
$ws.Cells.Item(383, 17).Value = 70
$ws.Cells.Item(383, 18).Value = "Hortaliza"
